$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows right before the old row 17 (BJT transistor row),
# pushing it (and everything below) down by three rows so we can slot the
# new LED components in above it.
$ws.Rows("17:19").Insert()

# Shop part numbers for the three new LEDs (typed first, matching the
# author's original entry order).
$ws.Range("C17").Value = "941-C4SMABGYCR34Q4S1"
$ws.Range("C18").Value = "941-C4SMAGGYCU2W37A1"
$ws.Range("C19").Value = "941-C4SMARGYCS4V1BB1"

# Component / description names.
$ws.Range("A17").Value = "Blue LED"
$ws.Range("B17").Value = "Blue LED"
$ws.Range("A18").Value = "Green LED"
$ws.Range("B18").Value = "Green LED"
$ws.Range("A19").Value = "Red LED"
$ws.Range("B19").Value = "Red LED"

# Prices, quantities and line totals.
$ws.Range("D17").Value = 0.106
$ws.Range("D18").Value = 0.101
$ws.Range("D19").Value = 0.101
$ws.Range("D17:D19").NumberFormat = "#,##0.00\ [$€-1];[Red]\-#,##0.00\ [$€-1]"

$ws.Range("E17").Value = 10
$ws.Range("E18").Value = 10
$ws.Range("E19").Value = 10

$ws.Range("F17").Formula = "=D17*E17"
$ws.Range("F18").Formula = "=D18*E18"
$ws.Range("F19").Formula = "=D19*E19"

# Shop part number (column C) and green/red price (column D) cells use a
# slightly darker grey font in the source workbook.
$ws.Range("C17").Font.Color = 3355443
$ws.Range("D18").Font.Color = 3355443
$ws.Range("D19").Font.Color = 3355443

# Printer/page setup, matching the author's resave in Croatian-localised Excel.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("F33").Select() | Out-Null
